$wb = $excel.ActiveWorkbook
$wsTests = $wb.Worksheets.Item("Tests")
$wsResult = $wb.Worksheets.Item("Result")

$data = @(
    @("Framework\InitAllSettings.xaml", "Success"),
    @("Framework\InitAllApplications.xaml", "Success"),
    @("Framework\CloseAllApplications.xaml", "Success"),
    @("Framework\CloseAllApplications.xaml", "SystemException"),
    @("Framework\InitAllSettings.xaml", "Success"),
    @("Framework\InitAllSettings.xaml", "Success"),
    @("Framework\InitAllApplications.xaml", "Success"),
    @("Framework\CloseAllApplications.xaml", "Success")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $wsResult.Cells.Item($row, 1).Value = $data[$i][0]
    $wsResult.Cells.Item($row, 2).Value = $data[$i][1]
}

$wsTests.Range("A2:B9").Select()

$wsResult.Activate()
$wsResult.Range("B5").Select()
